$wb = $excel.ActiveWorkbook

# --- Add "measure" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$measureSheet = $wb.Worksheets.Add($null, $lastSheet)
$measureSheet.Name = "measure"

$data7 = New-Object 'object[,]' 37,2
$data7[0,0] = 0
$data7[0,1] = -21.5
$data7[1,0] = 0
$data7[1,1] = -21.5
$data7[2,0] = 0.35355339059327401
$data7[2,1] = -19.1403499549969
$data7[3,0] = 0.42426406871192801
$data7[3,1] = -19.211060633115501
$data7[4,0] = 0.5
$data7[4,1] = -20
$data7[5,0] = 1.8
$data7[5,1] = -15
$data7[6,0] = 2
$data7[6,1] = -15
$data7[7,0] = 2.4041630560342599
$data7[7,1] = -14.119891808572399
$data7[8,0] = 2.4748737341529199
$data7[8,1] = -14.190602486691001
$data7[9,0] = 2.5
$data7[9,1] = -15
$data7[10,0] = 5
$data7[10,1] = -11.5
$data7[11,0] = 5
$data7[11,1] = -11.2
$data7[12,0] = 5.5154328932550696
$data7[12,1] = -10.1600938339277
$data7[13,0] = 5.6
$data7[13,1] = -10
$data7[14,0] = 5.6568542494923797
$data7[14,1] = -10.301515190165
$data7[15,0] = 6
$data7[15,1] = -10
$data7[16,0] = 9.3338095116624302
$data7[16,1] = -6.9074026404695799
$data7[17,0] = 9.5459415460183905
$data7[17,1] = -7.1195346748255499
$data7[18,0] = 10
$data7[18,1] = -7
$data7[19,0] = 10
$data7[19,1] = -6.7
$data7[20,0] = 13
$data7[20,1] = -5
$data7[21,0] = 13
$data7[21,1] = -5
$data7[22,0] = 13.930003589375
$data7[22,1] = -4.4325289063166604
$data7[23,0] = 14.142135623731001
$data7[23,1] = -4.6446609406726296
$data7[24,0] = 15
$data7[24,1] = -4.5
$data7[25,0] = 15
$data7[25,1] = -4.2
$data7[26,0] = 18.9504617357995
$data7[26,1] = -2.3819192408756802
$data7[27,0] = 19.091883092036799
$data7[27,1] = -2.5233405971129801
$data7[28,0] = 20
$data7[28,1] = -2.5
$data7[29,0] = 20
$data7[29,1] = -2.5
$data7[30,0] = 20.5
$data7[30,1] = -2.5
$data7[31,0] = 24.112341238461301
$data7[31,1] = -0.47273093167199698
$data7[32,0] = 24.395183950935898
$data7[32,1] = -0.75557364414661499
$data7[33,0] = 25
$data7[33,1] = -1.5
$data7[34,0] = 25
$data7[34,1] = -1.1000000000000001
$data7[35,0] = 29.5
$data7[35,1] = 0
$data7[36,0] = 29.5
$data7[36,1] = 0
$measureSheet.Range("A1:B37").Value = $data7

$measureSheet.Range("A1:B37").Select()

# --- Add "centered" sheet after "measure" ---
$centeredSheet = $wb.Worksheets.Add($null, $measureSheet)
$centeredSheet.Name = "centered"

$data8 = New-Object 'object[,]' 9,2
$data8[0,0] = 0
$data8[0,1] = -21.5
$data8[1,0] = 0.42593915310173402
$data8[1,1] = -19.4504701960375
$data8[2,0] = 2.23580735803744
$data8[2,1] = -14.662098859052699
$data8[3,0] = 5.4620478571245696
$data8[3,1] = -10.526934837348801
$data8[4,0] = 9.7199377644202105
$data8[4,1] = -6.9317343288237803
$data8[5,0] = 14.012023202184301
$data8[5,1] = -4.6295316411648804
$data8[6,0] = 19.708468965567299
$data8[6,1] = -2.4810519675977298
$data8[7,0] = 24.6268812973493
$data8[7,1] = -0.95707614395465301
$data8[8,0] = 29.5
$data8[8,1] = 0
$centeredSheet.Range("A1:B9").Value = $data8

$centeredSheet.Range("G32").Select()

# --- Restore selection on sheet "d" and re-activate it ---
$dSheet = $wb.Worksheets.Item("d")
$dSheet.Activate()
$dSheet.Range("H21").Select()
